# Update column G ("K") values on Sheet1 rows 2-26
# This regenerates the K column data (replacing former "Strike#" derived
# values) with newly calculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 3
    3  = 2
    4  = 2
    5  = 0
    6  = 3
    7  = 5
    8  = 2
    9  = 1
    10 = 5
    11 = 5
    12 = 4
    13 = 5
    14 = 3
    15 = 5
    16 = 4
    17 = 0
    18 = 1
    19 = 3
    20 = 2
    21 = 3
    22 = 3
    23 = 1
    24 = 1
    25 = 3
    26 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
